$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 25
$ws1.Range("F3").Value = 16474
$ws1.Range("F5").Value = 741
$ws1.Range("F7").Value = 73
$ws1.Range("F8").Value = 9291
$ws1.Range("F11").Value = 1036
$ws1.Range("F12").Value = 131
$ws1.Range("F13").Value = 228
$ws1.Range("F18").Value = 631
$ws1.Range("F20").Value = 17
$ws1.Range("F21").Value = 80
$ws1.Range("F26").Value = 541
$ws1.Range("F27").Value = 41
$ws1.Range("F32").Value = 69
$ws1.Range("F33").Value = 274
$ws1.Range("F34").Value = 376
$ws1.Range("F37").Value = 5741
$ws1.Range("F38").Value = 5257

# Sheet "演出" (Show) - column F update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81

# Sheet "全部类型" (All types) - column F "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 25
$ws4.Range("F3").Value = 16474
$ws4.Range("F5").Value = 741
$ws4.Range("F7").Value = 73
$ws4.Range("F8").Value = 9291
$ws4.Range("F11").Value = 1036
$ws4.Range("F12").Value = 131
$ws4.Range("F13").Value = 228
$ws4.Range("F18").Value = 631
$ws4.Range("F20").Value = 17
$ws4.Range("F21").Value = 80
$ws4.Range("F26").Value = 541
$ws4.Range("F27").Value = 41
$ws4.Range("F29").Value = 81
$ws4.Range("F34").Value = 69
$ws4.Range("F35").Value = 274
$ws4.Range("F36").Value = 376
$ws4.Range("F39").Value = 5741
$ws4.Range("F41").Value = 5257
